$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NOVIEMBRE 2024")

# ---------------------------------------------------------------------------
# 1) The last two "open" days (45870, 45871) that had no activity logged yet
#    are removed now that the week is being closed out and a new week's
#    worth of entries is being recorded below.
# ---------------------------------------------------------------------------
$ws.Rows("265:266").Delete()

# ---------------------------------------------------------------------------
# 2) Build the new "SEMANA 30" block by duplicating the structure of the
#    previous week's block (header + column titles + data rows), then
#    overwrite the copied text with the new week's content. Using separate,
#    non-overlapping copy operations keeps formatting/merges identical to
#    the rest of the sheet.
# ---------------------------------------------------------------------------
$ws.Range("B261:F264").Copy($ws.Range("B267"))
$ws.Range("B263:F264").Copy($ws.Range("B271"))
$ws.Range("B263:F264").Copy($ws.Range("B273"))

# Week header
$ws.Range("B267").Value = "SEMANA 30"

# Row 269 - 2025-08-04
$ws.Range("B269").Value2 = 45873
$ws.Range("C269").Value = "Se agregarón efectos con observadores"
$ws.Range("D269").Value = "Se agregarón efectos con observadores"
$ws.Range("E269").Value = "8:00 - 13:00, 17:00 - 18:00"
$ws.Range("F269").Value = 6

# Row 270 - 2025-08-04
$ws.Range("B270").Value2 = 45873
$ws.Range("C270").Value = "Se trabajó actividades de la OGA"
$ws.Range("D270").Value = "Se realizarón las actividades solicitadas de la OGA"
$ws.Range("E270").Value = "10:00 - 13:00, 15:00 -18:00"
$ws.Range("F270").Value = 6

# Row 271 - 2025-08-05 (no activity)
$ws.Range("B271").Value2 = 45874
$ws.Range("C271").Value = "-"
$ws.Range("D271").Value = "-"
$ws.Range("E271").Value = "-"
$ws.Range("F271").Value = "-"

# Row 272 - 2025-08-06 (no activity)
$ws.Range("B272").Value2 = 45875
$ws.Range("C272").Value = "-"
$ws.Range("D272").Value = "-"
$ws.Range("E272").Value = "-"
$ws.Range("F272").Value = "-"

# Row 273 - 2025-08-07
$ws.Range("B273").Value2 = 45876
$ws.Range("C273").Value = "Se trabajo en los cambios solicitados en la última reuinion"
$ws.Range("D273").Value = "Se trabajo en los cambios solicitados en la última reunión, para la planta de tratamiento AR."
$ws.Range("E273").Value = "8:00 -12:00, 15:00 - 18:00"
$ws.Range("F273").Value = 7

# Row 274 - 2025-08-08 (not filled in yet)
$ws.Range("B274").Value2 = 45877
$ws.Range("C274").ClearContents()
$ws.Range("D274").ClearContents()
$ws.Range("E274").ClearContents()
$ws.Range("F274").ClearContents()

# Weekly total (copy the style from the previous week's total cell, then
# give it this week's formula)
$ws.Range("G265").Copy($ws.Range("G275"))
$ws.Range("G275").Formula = "=SUM(F269:F274)"
